$wb = $excel.ActiveWorkbook

# --- Update Sheet13: add new detection box data (row 7, rows 18-19) ---
$sheet13 = $wb.Worksheets.Item("Sheet13")

$sheet13.Range("A7").Value = 54
$sheet13.Range("B7").Value = 36
$sheet13.Range("E7").Value = 72
$sheet13.Range("G7").Value = 90

$sheet13.Range("C18").Value = 54
$sheet13.Range("I18").Value = 0
$sheet13.Range("E19").Value = 36

$sheet13.Activate()
$sheet13.Range("G7").Select()

# --- Add Sheet14 (new detection results) after Sheet13 ---
$sheet14 = $wb.Worksheets.Add($null, $sheet13)
$sheet14.Name = "Sheet14"

$sheet14.Range("B1").Value = 96

$sheet14.Range("C4").Value = 90
$sheet14.Range("E4").Value = 108
$sheet14.Range("G4").Value = 126

$sheet14.Range("C7").Value = 0
$sheet14.Range("E7").Value = 0
$sheet14.Range("G7").Value = 18
$sheet14.Range("I7").Value = 36

$sheet14.Range("C10").Value = 36
$sheet14.Range("D10").Value = 36

$sheet14.Range("D14").Value = 36
$sheet14.Range("E14").Value = 36

$sheet14.Range("B16").Value = 67

$sheet14.Range("H18").Value = 67

$sheet14.Range("B19").Value = 67

$sheet14.Range("G21").Value = 96

$sheet14.Range("B22").Value = 67
$sheet14.Range("C22").Value = 67

$sheet14.Range("J23").Value = 96

$sheet14.Range("B1").Select()

# --- Add Sheet15 (new detection results) after Sheet14 ---
$sheet15 = $wb.Worksheets.Add($null, $sheet14)
$sheet15.Name = "Sheet15"

$sheet15.Range("B3").Value = 0
$sheet15.Range("F3").Value = 0

$sheet15.Range("B7").Value = 0
$sheet15.Range("F7").Value = 0

$sheet15.Range("F8").Value = 0

$sheet15.Range("B13").Value = 0
$sheet15.Range("C13").Value = 0
$sheet15.Range("D13").Value = 0
$sheet15.Range("E13").Value = 0

$sheet15.Range("B14").Value = 36
$sheet15.Range("C14").Value = 54

$sheet15.Range("C15").Value = 0

# Sheet15 is the last/active sheet and has the selection on C15
$sheet15.Activate()
$sheet15.Range("C15").Select()
